# Changes to the Shop Order Operation table
# - Insert three new columns into the "ShopOrderOperations" sheet:
#     F: PrecedingOperation
#     G: WCRuntimeFactor
#     I: LaborRuntimeFactor
#   (existing WorkCenterRuntime / LaborRuntime / OpStartDate.. etc shift right)
# - Populate the new columns with data
# - Make "ShopOrderOperations" the active sheet/tab, with F11 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShopOrderOperations")

# Insert the 3 new columns (order matters: F, then G, then I)
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()

# New header labels
$ws.Range("F1").Value = "PrecedingOperation"
$ws.Range("G1").Value = "WCRuntimeFactor"
$ws.Range("I1").Value = "LaborRuntimeFactor"

# New column data values
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 0

$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 0

$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("I4").Value = 0

$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 0

# Make ShopOrderOperations the active sheet and select F11
$ws.Activate()
$ws.Range("F11").Select()
